$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column: header in H1 (matching the formatting of the other
# header cells like G1) and its data value in H2.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
